# Updated cryptos list: refreshed Price/Volume(1h) figures for every coin
# row, plus TheGraph and ApeXProtocol swapping rank positions (44 <-> 45).
# "D" column Price values are text (not numbers) in the source data, so
# they're written with a leading apostrophe to force Excel's text storage
# instead of letting it auto-parse them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.933.45"
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = "'3.316.95"
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = "'582.39"
$ws.Range("E5").Value = '  -1.41%  '
$ws.Range("D6").Value = "'175.64"
$ws.Range("E6").Value = '  -6.18%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("D9").Value = "'3.308.98"
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("E10").Value = '  -4.00%  '
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("D12").Value = "'45.46"
$ws.Range("E12").Value = '  -4.32%  '
$ws.Range("E13").Value = '  -2.17%  '
$ws.Range("D14").Value = "'665.17"
$ws.Range("E14").Value = '  +5.00%  '
$ws.Range("D15").Value = "'3.853.99"
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("E16").Value = '  -2.75%  '
$ws.Range("D17").Value = "'67.926.40"
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").Value = "'3.312.41"
$ws.Range("E19").Value = '  -1.69%  '
$ws.Range("E20").Value = '  -3.07%  '
$ws.Range("D21").Value = "'10.88"
$ws.Range("E22").Value = '  -2.37%  '
$ws.Range("E23").Value = '  +5.30%  '
$ws.Range("D24").Value = "'17.18"
$ws.Range("E24").Value = '  -4.55%  '
$ws.Range("D25").Value = "'97.60"
$ws.Range("E25").Value = '  -1.85%  '
$ws.Range("E26").Value = '  -4.20%  '
$ws.Range("E27").Value = '  -5.87%  '
$ws.Range("D28").Value = "'9.25"
$ws.Range("E28").Value = '  -4.49%  '
$ws.Range("D29").Value = "'33.41"
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("E30").Value = '  -2.98%  '
$ws.Range("E31").Value = '  +2.63%  '
$ws.Range("D32").Value = "'587.13"
$ws.Range("E32").Value = '  -2.92%  '
$ws.Range("D33").Value = "'10.95"
$ws.Range("E33").Value = '  -1.33%  '
$ws.Range("E34").Value = '  -1.75%  '
$ws.Range("D35").Value = "'3.750.66"
$ws.Range("E35").Value = '  -5.66%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").Value = "'3.35"
$ws.Range("E37").Value = '  -12.10%  '
$ws.Range("D38").Value = "'55.36"
$ws.Range("E38").Value = '  -1.25%  '
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("D40").Value = "'2.62"
$ws.Range("E40").Value = '  -7.21%  '
$ws.Range("D41").Value = "'32.37"
$ws.Range("E41").Value = '  -4.22%  '
$ws.Range("D42").Value = "'3.14"
$ws.Range("E42").Value = '  -3.11%  '
$ws.Range("E43").Value = '  -5.27%  '

# Row 44/45: TheGraph and ApeXProtocol swap rank positions (TheGraph now
# ranked 44th, ApeXProtocol 45th), each with updated Price/Volume figures.
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.332"
$ws.Range("E44").Value = "  -3.08%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.26"
$ws.Range("E45").Value = "  -4.10%  "

$ws.Range("E46").Value = '  -3.95%  '
$ws.Range("D47").Value = "'2.60"
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("E48").Value = '  -2.23%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = "'1.34"
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("D51").Value = "'129.57"
$ws.Range("E51").Value = '  +0.04%  '